$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3 corresponds to file_name = metrics_sim_with_priors.json
# Updating recall / wss / loss / erf / atd metrics
$ws.Range("C3").Value = 0.3303571428571428
$ws.Range("D3").Value = 0.6041666666666666
$ws.Range("F3").Value = 0.9910714285714286
$ws.Range("H3").Value = 0.3390952556086796
$ws.Range("I3").Value = 0.1920248586215854
$ws.Range("J3").Value = 0.2321428571428572
$ws.Range("K3").Value = 626.0952380952381

# Updating fp_* counts
$ws.Range("Q3").Value = 67
$ws.Range("R3").Value = 109
$ws.Range("S3").Value = 324
$ws.Range("T3").Value = 691
$ws.Range("U3").Value = 1071

# Updating tn_* counts
$ws.Range("V3").Value = 2316
$ws.Range("W3").Value = 2274
$ws.Range("X3").Value = 2059
$ws.Range("Y3").Value = 1692
$ws.Range("Z3").Value = 1312

# Updating tnr_* metrics
$ws.Range("AF3").Value = 0.971884
$ws.Range("AG3").Value = 0.954259
$ws.Range("AH3").Value = 0.8640370000000001
$ws.Range("AI3").Value = 0.710029
$ws.Range("AJ3").Value = 0.550567
